$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header column C from "audioFalse" to "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# Collapse the two distinct audio-file values in column C (which were
# representing the "audioFalse" condition) into a single "train1P2" value
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
